$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp footer (last row, column A)
$ws.Range("A1").Value = "Datos actualizados a 14 de Junio de 2020 a las 21:31"

# Apply country data updates (values refreshed + four pairs of rows swapped
# position due to re-sorting by "Casos totales")
$changes = @(
    @{Addr="B4"; Value=2156772},
    @{Addr="C4"; Value=14548},
    @{Addr="D4"; Value=857309},
    @{Addr="E4"; Value=1181735},
    @{Addr="G4"; Value=201},
    @{Addr="H4"; Value=117728},
    @{Addr="B7"; Value=332970},
    @{Addr="C7"; Value=11344},
    @{Addr="D7"; Value=169600},
    @{Addr="E7"; Value=153850},
    @{Addr="G7"; Value=321},
    @{Addr="H7"; Value=9520},
    @{Addr="B12"; Value=187631},
    @{Addr="C12"; Value=208},
    @{Addr="E12"; Value=6562},
    @{Addr="G12"; Value=2},
    @{Addr="H12"; Value=8869},
    @{Addr="D15"; Value=143704},
    @{Addr="E15"; Value=27266},
    @{Addr="B76"; Value=5080},
    @{Addr="C76"; Value=114},
    @{Addr="E76"; Value=1118},
    @{Addr="A107"; Value="Costa Rica"},
    @{Addr="B107"; Value=1715},
    @{Addr="C107"; Value=53},
    @{Addr="D107"; Value=752},
    @{Addr="E107"; Value=951},
    @{Addr="H107"; Value=12},
    @{Addr="A108"; Value="Sudan del Sur"},
    @{Addr="B108"; Value=1693},
    @{Addr="D108"; Value=49},
    @{Addr="E108"; Value=1617},
    @{Addr="H108"; Value=27},
    @{Addr="A109"; Value="Mauritania"},
    @{Addr="B109"; Value=1682},
    @{Addr="D109"; Value=311},
    @{Addr="E109"; Value=1288},
    @{Addr="H109"; Value=83},
    @{Addr="A144"; Value="Ruanda"},
    @{Addr="B144"; Value=582},
    @{Addr="C144"; Value=41},
    @{Addr="D144"; Value=332},
    @{Addr="E144"; Value=248},
    @{Addr="G144"; Value=0},
    @{Addr="H144"; Value=2},
    @{Addr="A145"; Value="Malaui"},
    @{Addr="B145"; Value=547},
    @{Addr="C145"; Value=18},
    @{Addr="D145"; Value=69},
    @{Addr="E145"; Value=472},
    @{Addr="G145"; Value=1},
    @{Addr="H145"; Value=6},
    @{Addr="A153"; Value="Benin"},
    @{Addr="B153"; Value=442},
    @{Addr="C153"; Value=30},
    @{Addr="D153"; Value=228},
    @{Addr="E153"; Value=208},
    @{Addr="H153"; Value=6},
    @{Addr="A154"; Value="Libia"},
    @{Addr="B154"; Value=418},
    @{Addr="D154"; Value=62},
    @{Addr="E154"; Value=348},
    @{Addr="H154"; Value=8},
    @{Addr="B155"; Value=383},
    @{Addr="C155"; Value=27},
    @{Addr="E155"; Value=325},
    @{Addr="B175"; Value=118},
    @{Addr="C175"; Value=1},
    @{Addr="E175"; Value=1},
    @{Addr="A206"; Value="Groenlandia"},
    @{Addr="A207"; Value="Islas Malvinas"},
    @{Addr="A208"; Value="Islas Turcas y Caicos"},
    @{Addr="D208"; Value=11},
    @{Addr="H208"; Value=1},
    @{Addr="A209"; Value="Santa Sede"},
    @{Addr="D209"; Value=12},
    @{Addr="H209"; Value=0}
)

foreach ($ch in $changes) {
    $ws.Range($ch.Addr).Value = $ch.Value
}

